# Applies the row-data changes from the diff: rows 2-8 of the "Artfynd"
# worksheet are updated in place (values rotate among the existing rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 94980872
$ws.Cells.Item(2, 2).Value = 56411
$ws.Cells.Item(2, 4).Value = "'NT"
$ws.Cells.Item(2, 5).Value = 100049
$ws.Cells.Item(2, 6).Value = "'Spillkråka"
$ws.Cells.Item(2, 7).Value = "'Dryocopus martius"
$ws.Cells.Item(2, 8).Value = "'(Linnaeus, 1758)"
$ws.Cells.Item(2, 9).Value = ""
$ws.Cells.Item(2, 11).Value = ""
$ws.Cells.Item(2, 13).Value = "'förbiflygande"
$ws.Cells.Item(2, 17).Value = 572435.745707145
$ws.Cells.Item(2, 18).Value = 6701433.942975808
$ws.Cells.Item(2, 26).Value = "'00:00"
$ws.Cells.Item(2, 28).Value = "'00:00"

# Row 3
$ws.Cells.Item(3, 1).Value = 94979072
$ws.Cells.Item(3, 2).Value = 56411
$ws.Cells.Item(3, 4).Value = "'NT"
$ws.Cells.Item(3, 5).Value = 100049
$ws.Cells.Item(3, 6).Value = "'Spillkråka"
$ws.Cells.Item(3, 7).Value = "'Dryocopus martius"
$ws.Cells.Item(3, 8).Value = "'(Linnaeus, 1758)"
$ws.Cells.Item(3, 9).Value = ""
$ws.Cells.Item(3, 11).Value = ""
$ws.Cells.Item(3, 13).Value = "'äldre spår"
$ws.Cells.Item(3, 17).Value = 572416.4402874345
$ws.Cells.Item(3, 18).Value = 6701459.721651074
$ws.Cells.Item(3, 26).Value = "'20:24"
$ws.Cells.Item(3, 28).Value = "'20:24"

# Row 4
$ws.Cells.Item(4, 1).Value = 94979361
$ws.Cells.Item(4, 9).Value = "'24"
$ws.Cells.Item(4, 16).Value = "'Skurtjärn, Dlr"
$ws.Cells.Item(4, 17).Value = 572430.1908244109
$ws.Cells.Item(4, 18).Value = 6701241.788604234
$ws.Cells.Item(4, 26).Value = "'20:39"
$ws.Cells.Item(4, 28).Value = "'20:39"
$ws.Cells.Item(4, 49).Value = "'Philipp Weiss"
$ws.Cells.Item(4, 50).Value = "'Philipp Weiss"

# Row 5
$ws.Cells.Item(5, 1).Value = 94979464
$ws.Cells.Item(5, 9).Value = "'6"
$ws.Cells.Item(5, 17).Value = 572409.0058300388
$ws.Cells.Item(5, 18).Value = 6701262.594218019
$ws.Cells.Item(5, 26).Value = "'20:43"
$ws.Cells.Item(5, 28).Value = "'20:43"

# Row 6
$ws.Cells.Item(6, 1).Value = 94979470
$ws.Cells.Item(6, 2).Value = 96334
$ws.Cells.Item(6, 4).Value = "'VU"
$ws.Cells.Item(6, 5).Value = 220787
$ws.Cells.Item(6, 6).Value = "'Knärot"
$ws.Cells.Item(6, 7).Value = "'Goodyera repens"
$ws.Cells.Item(6, 8).Value = "'(L.) R. Br."
$ws.Cells.Item(6, 11).Value = "'blomning"
$ws.Cells.Item(6, 13).Value = ""
$ws.Cells.Item(6, 16).Value = "'Hedemora, Dlr"
$ws.Cells.Item(6, 17).Value = 572458.1188458267
$ws.Cells.Item(6, 18).Value = 6701230.004390508
$ws.Cells.Item(6, 26).Value = "'20:43"
$ws.Cells.Item(6, 28).Value = "'20:43"
$ws.Cells.Item(6, 49).Value = "'Mariapaz Ojeda"
$ws.Cells.Item(6, 50).Value = "'Mariapaz Ojeda"

# Row 7
$ws.Cells.Item(7, 1).Value = 94979424
$ws.Cells.Item(7, 2).Value = 96334
$ws.Cells.Item(7, 4).Value = "'VU"
$ws.Cells.Item(7, 5).Value = 220787
$ws.Cells.Item(7, 6).Value = "'Knärot"
$ws.Cells.Item(7, 7).Value = "'Goodyera repens"
$ws.Cells.Item(7, 8).Value = "'(L.) R. Br."
$ws.Cells.Item(7, 11).Value = "'blomning"
$ws.Cells.Item(7, 13).Value = ""
$ws.Cells.Item(7, 17).Value = 572426.1968956392
$ws.Cells.Item(7, 18).Value = 6701243.683577545
$ws.Cells.Item(7, 26).Value = "'20:41"
$ws.Cells.Item(7, 28).Value = "'20:41"

# Row 8
$ws.Cells.Item(8, 1).Value = 94979445
$ws.Cells.Item(8, 2).Value = 5113
$ws.Cells.Item(8, 4).Value = "'LC"
$ws.Cells.Item(8, 5).Value = 100526
$ws.Cells.Item(8, 6).Value = "'Bronshjon"
$ws.Cells.Item(8, 7).Value = "'Callidium coriaceum"
$ws.Cells.Item(8, 8).Value = "'Paykull, 1800"
$ws.Cells.Item(8, 13).Value = "'äldre gnagspår"
$ws.Cells.Item(8, 17).Value = 572423.8744920741
$ws.Cells.Item(8, 18).Value = 6701260.916536125
$ws.Cells.Item(8, 26).Value = "'20:42"
$ws.Cells.Item(8, 28).Value = "'20:42"
